$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 131126459
$ws.Range("B3").Value = 58043
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 103021
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I3").Value = "'1"
$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("P3").Value = "Mölleröd, Sk"
$ws.Range("Q3").Value = 452800
$ws.Range("R3").Value = 6221182
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Skåne"
$ws.Range("U3").Value = "Kristianstad"
$ws.Range("V3").Value = "Skåne"
$ws.Range("W3").Value = "Österslöv"
$ws.Range("Y3").Value = "'2026-02-12"
$ws.Range("AA3").Value = "'2026-02-12"
$ws.Range("AC3").Value = "I meståg."
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AT3").Value = "'"
$ws.Range("AW3").Value = "Roine Strandberg"
$ws.Range("AX3").Value = "Roine Strandberg"
$ws.Range("AY3").Value = "'"

$ws.Range("I3").Style = "Normal"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").Style = "Normal"
$ws.Range("N3").Style = "Normal"
$ws.Range("Y3").Style = "Normal"
$ws.Range("AA3").Style = "Normal"
$ws.Range("AT3").Style = "Normal"
$ws.Range("AY3").Style = "Normal"
